$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 / Row 12 swap: TRON and Dogecoin trade places ---
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.134"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0845"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.19%  "

# --- Remaining Price (D) / Volume 1h (E) updates ---
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.861.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.781.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.07%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "356.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.558"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.48%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.56%  "

# Row 13
$ws.Range("E13").Value = "  -3.11%  "

# Row 14
$ws.Range("E14").Value = "  -3.36%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.228.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.77%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.788.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.49%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.938"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.91%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.859.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.92%  "

# Row 20
$ws.Range("E20").Value = "  -2.41%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.63%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.52%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.46%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.14%  "

# Row 25
$ws.Range("E25").Value = "  -3.46%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.47%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +16.55%  "

# Row 29
$ws.Range("E29").Value = "  -1.23%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.96%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0469"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.84%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.46%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.37%  "

# Row 34
$ws.Range("E34").Value = "  -2.73%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0843"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.31%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.03%  "

# Row 37
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.18%  "

# Row 40
$ws.Range("E40").Value = "  -4.37%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.20%  "

# Row 42
$ws.Range("E42").Value = "  -2.52%  "

# Row 43
$ws.Range("E43").Value = "  -1.05%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.63%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.04%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.091.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.10%  "

# Row 47
$ws.Range("E47").Value = "  -5.09%  "

# Row 48
$ws.Range("E48").Value = "  -1.79%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.62%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.952"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.85%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.07%  "
